$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.3103
$ws.Range("C6").Value = -12.81709999999999
$ws.Range("C7").Value = -12.44309999999999
$ws.Range("D7").Value = -7.515799999999993
$ws.Range("C8").Value = -13.8815
$ws.Range("D11").Value = -7.743699999999999
$ws.Range("D12").Value = -7.161799999999997
$ws.Range("E12").Value = 17.3232
$ws.Range("E13").Value = 16.89730000000001
$ws.Range("E14").Value = 17.077
$ws.Range("D15").Value = -9.005399999999989
$ws.Range("C16").Value = -14.01409999999999
$ws.Range("E16").Value = 16.25370000000001
$ws.Range("E19").Value = 16.44380000000001
$ws.Range("C20").Value = -12.0569
$ws.Range("D20").Value = -8.273299999999995
$ws.Range("E20").Value = 16.8066
$ws.Range("C21").Value = -12.6737
$ws.Range("D21").Value = -8.107799999999999
$ws.Range("D22").Value = -8.091000000000003
$ws.Range("E22").Value = 16.4845
$ws.Range("D23").Value = -7.405399999999998
$ws.Range("C28").Value = -12.494
$ws.Range("C29").Value = -11.12760000000001
$ws.Range("D29").Value = -7.079899999999998
$ws.Range("C30").Value = -13.16359999999999
$ws.Range("C32").Value = -12.5016
$ws.Range("D34").Value = -7.863400000000002
$ws.Range("E36").Value = 15.6658
$ws.Range("C40").Value = -11.76060000000001
$ws.Range("D42").Value = -8.737699999999997
$ws.Range("D43").Value = -8.101699999999999
$ws.Range("E43").Value = 16.9574
$ws.Range("D44").Value = -7.122199999999997
$ws.Range("D45").Value = -7.619599999999998
$ws.Range("C46").Value = -14.69019999999999
$ws.Range("D46").Value = -8.335900000000004
$ws.Range("E46").Value = 16.54840000000002
$ws.Range("D50").Value = -7.971399999999996
$ws.Range("E50").Value = 16.96139999999999
$ws.Range("C51").Value = -12.392
$ws.Range("D51").Value = -7.724599999999997
$ws.Range("C52").Value = -11.17510000000001
$ws.Range("C57").Value = -13.91079999999999
$ws.Range("D57").Value = -8.754199999999999
$ws.Range("C59").Value = -12.77779999999999
$ws.Range("C62").Value = -14.10589999999998
$ws.Range("D65").Value = -7.768599999999997
$ws.Range("C66").Value = -13.2203
$ws.Range("D66").Value = -7.806899999999998
$ws.Range("D67").Value = -6.422100000000004
$ws.Range("C73").Value = -11.26
$ws.Range("C74").Value = -12.18120000000001
$ws.Range("E76").Value = 16.41699999999999
$ws.Range("C77").Value = -11.6395
$ws.Range("D79").Value = -6.288100000000006
$ws.Range("D84").Value = -8.900100000000002
$ws.Range("D87").Value = -7.981999999999999
$ws.Range("C92").Value = -11.94320000000001
$ws.Range("D92").Value = -6.672700000000003
$ws.Range("E95").Value = 18.43480000000003
$ws.Range("D97").Value = -8.552099999999999
$ws.Range("E97").Value = 16.22899999999999
$ws.Range("E99").Value = 16.584
$ws.Range("C100").Value = -12.3055
